$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 was Facebook/Likes -> becomes Twitter/Follower with new ThisWeek/LastWeek values
$ws.Range("A3").Value = "Twitter"
$ws.Range("B3").Value = "Follower"
$ws.Range("C3").Value = 4654
$ws.Range("D3").Value = 4684

# Row 4 was Twitter/Follower -> becomes Instagram/Follower with the old row 5's values
$ws.Range("A4").Value = "Instagram"
$ws.Range("B4").Value = "Follower"
$ws.Range("C4").Value = 5416
$ws.Range("D4").Value = 4683

# Old row 5 (Instagram) is no longer needed - remove it entirely
$ws.Rows("5").Delete()

# Update the active selection to match the saved workbook state
$ws.Range("D3").Select()
